$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (r2_adj) - the table shrinks from 5 rows to 4 rows.
$ws.Rows("5:5").Delete()

# Add the new "C/A" column by copying the formatting of the existing
# header cell (C1) into D1 before filling in the new header/data text -
# this reuses the bold/bordered/centered header style instead of
# introducing a new one.
$ws.Range("C1").Copy($ws.Range("D1"))

# Header row: FFR | C/A | $\pi$
$ws.Range("B1").Value = "FFR"
$ws.Range("C1").Value = "C/A"
$ws.Range("D1").Value = '$\pi$'

# Row labels (A2:A4) - order is now FFR Lag, C/A Lag, $\pi$ Lag
$ws.Range("A2").Value = "FFR Lag"
$ws.Range("A3").Value = "C/A Lag"
$ws.Range("A4").Value = '$\pi$ Lag'

# Data values, row 2 (FFR Lag)
$ws.Range("B2").Value = "0.424***"
$ws.Range("C2").Value = "9.463***"
$ws.Range("D2").Value = "0.352***"

# Data values, row 3 (C/A Lag)
$ws.Range("B3").Value = "-0.033***"
$ws.Range("C3").Value = "-0.763***"
$ws.Range("D3").Value = "-0.016***"

# Data values, row 4 ($\pi$ Lag). B4's "0.039" is purely numeric-looking,
# so force it to stay text (matching the other text-typed coefficients)
# by temporarily applying a text number format, then reverting the style
# so the cell itself keeps its original (unstyled) appearance.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "0.039"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "2.634***"
$ws.Range("D4").Value = "-0.656***"
